$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12704
$ws1.Range("F3").Value = 615
$ws1.Range("F5").Value = 21
$ws1.Range("F6").Value = 290
$ws1.Range("F7").Value = 397
$ws1.Range("F9").Value = 12701
$ws1.Range("F10").Value = 27
$ws1.Range("F11").Value = 3298
$ws1.Range("F12").Value = 543
$ws1.Range("F13").Value = 12
$ws1.Range("F14").Value = 9
$ws1.Range("F15").Value = 22
$ws1.Range("F16").Value = 1198
$ws1.Range("F17").Value = 27
$ws1.Range("F18").Value = 131
$ws1.Range("F19").Value = 662
$ws1.Range("F20").Value = 2844
$ws1.Range("F21").Value = 6129
$ws1.Range("F22").Value = 296
$ws1.Range("F23").Value = 3616

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12704
$ws4.Range("F3").Value = 615
$ws4.Range("F5").Value = 21
$ws4.Range("F6").Value = 290
$ws4.Range("F8").Value = 397
$ws4.Range("F10").Value = 12701
$ws4.Range("F11").Value = 27
$ws4.Range("F12").Value = 3299
$ws4.Range("F13").Value = 543
$ws4.Range("F14").Value = 12
$ws4.Range("F15").Value = 9
$ws4.Range("F16").Value = 22
$ws4.Range("F17").Value = 1198
$ws4.Range("F18").Value = 27
$ws4.Range("F19").Value = 131
$ws4.Range("F20").Value = 662
$ws4.Range("F21").Value = 2844
$ws4.Range("F23").Value = 6129
$ws4.Range("F24").Value = 297
$ws4.Range("F25").Value = 3616
